# Generate Report for Handoff
# A new source file (f6793463-...) replaces the old one (7d71bcde-...) and a
# fresh set of handoff .xlf files is generated. Because this report reflects
# a just-generated handoff, the "target"/"handback" columns for each locale
# are reset (no translation has come back yet).

$wb = $excel.ActiveWorkbook

$oldId  = "7d71bcde-2188-4dea-9990-360038628121"
$newId  = "f6793463-7e8f-4352-8c3e-76bc8cd562d6"
$oldHash = "d4ae73b8332663ba5dcae6116e1de887174e2bbf"
$newHash = "977fe874e62be926dccac48336d35741c55b0292"

$newFileName      = "$newId.md"
$newPathAndName   = "e2e\$newId.md"
$newGenerateDate  = "2016-09-02 07:08:58"

$newZhHandoffFile = "$newId.$newHash.zh-cn.xlf"
$newZhHandoffDate = "2016-09-02 07:08:54"
$newDeHandoffFile = "$newId.$newHash.de-de.xlf"

$resetDateTime = "0001-01-01 00:00:00"

function Set-HyperlinkDisplay($ws, $addr, $text) {
    foreach ($hl in $ws.Hyperlinks) {
        if ($hl.Range.Address() -eq $addr) {
            $hl.TextToDisplay = $text
        }
    }
}

function Remove-HyperlinkAt($ws, $addr) {
    foreach ($hl in $ws.Hyperlinks) {
        if ($hl.Range.Address() -eq $addr) {
            $hl.Delete()
        }
    }
}

# ---------------------------------------------------------------------------
# Overview sheet
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("A2").Value = $newFileName
Set-HyperlinkDisplay $wsOverview '$B$2' $newPathAndName
$wsOverview.Range("G2").Value = $newGenerateDate

# ---------------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

$wsZh.Range("A2").Value = $newFileName
Set-HyperlinkDisplay $wsZh '$A$2' $newFileName

$wsZh.Range("G2").Value = $newZhHandoffFile
$wsZh.Range("H2").Value = $newZhHandoffDate

Remove-HyperlinkAt $wsZh '$I$2'
$wsZh.Range("I2").ClearContents()
$wsZh.Range("I2").ClearFormats()

$wsZh.Range("J2").Value = ""
$wsZh.Range("K2").Value = $resetDateTime

# Columns I/J auto-shrink now that their (previously 40-char-wide) UUID
# content is gone; re-fit them to the remaining header text.
$wsZh.Columns.Item(9).ColumnWidth = 17.833333333333332
$wsZh.Columns.Item(10).ColumnWidth = 20.833333333333332

# ---------------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

$wsDe.Range("A2").Value = $newFileName
Set-HyperlinkDisplay $wsDe '$A$2' $newFileName

$wsDe.Range("G2").Value = $newDeHandoffFile
$wsDe.Range("H2").Value = $newGenerateDate

Remove-HyperlinkAt $wsDe '$I$2'
$wsDe.Range("I2").ClearContents()
$wsDe.Range("I2").ClearFormats()

$wsDe.Range("J2").Value = ""
$wsDe.Range("K2").Value = $resetDateTime

$wsDe.Columns.Item(9).ColumnWidth = 17.833333333333332
$wsDe.Columns.Item(10).ColumnWidth = 20.833333333333332
